# Add a new client record (C1003 / "hola") as row 5 of the Clientes sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Clientes")

# Columns E (fecha_ingreso) and F (fecha_dispersion) hold date-looking text
# ("2025-10-08") that must stay plain text, like the existing rows 3/4 —
# format the cells as Text first so the engine doesn't auto-convert the
# value to a serial date, then drop back to the default "Normal" style so
# no stray formatting is left behind on the cells.
$ws.Range("E5:F5").NumberFormat = "@"

$ws.Range("A5").Value = "C1003"
$ws.Range("B5").Value = "hola"
$ws.Range("C5").Value = "TOXQUI"
$ws.Range("D5").Value = "Martha Ortiz"
$ws.Range("E5").Value = "2025-10-08"
$ws.Range("F5").Value = "2025-10-08"
$ws.Range("G5").Value = "DISPERSADO"

$ws.Range("E5:F5").Style = "Normal"
